$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chris")
$ws.Activate()

# Excel color longs (BGR order): yellow fill (FFFF00) and blue fill (0070C0)
$yellow = 65535
$blue   = 12611584
$xlHAlignRight = -4152

# Rows marked as "addressed" with the yellow highlight (matches existing B-column yellow rows)
$yellowRows = @(4, 27, 28, 42, 47, 49, 61, 75)
foreach ($r in $yellowRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Interior.Color = $yellow
    $cell.HorizontalAlignment = $xlHAlignRight
}

# Rows marked as "addressed" with the blue highlight (matches existing B-column blue rows)
$blueRows = @(56, 57, 58)
foreach ($r in $blueRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Interior.Color = $blue
    $cell.HorizontalAlignment = $xlHAlignRight
}

# Rows marked with an "x" comment-resolution marker (no fill)
$xRows = @(5, 6, 7, 29)
foreach ($r in $xRows) {
    $ws.Cells.Item($r, 1).Value = "x"
}

# Update the view: scroll position and current selection
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("A29").Select()
